$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.975.14"
$ws.Range("E2").Value = "  -2.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.021.77"
$ws.Range("E3").Value = "  -2.06%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.00"
$ws.Range("E5").Value = "  +1.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.58"
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.016.13"
$ws.Range("E8").Value = "  -2.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("E10").Value = "  -4.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.10"
$ws.Range("E11").Value = "  -5.00%  "
$ws.Range("E12").Value = "  -1.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000223"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.35"
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.513.14"
$ws.Range("E15").Value = "  -2.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.060.97"
$ws.Range("E16").Value = "  -2.79%  "
$ws.Range("E17").Value = "  -2.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.027.14"
$ws.Range("E18").Value = "  -1.96%  "
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "475.51"
$ws.Range("E20").Value = "  -1.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.28"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("E22").Value = "  -3.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.07"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.45"
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.14"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.75"
$ws.Range("E28").Value = "  -3.80%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.76"
$ws.Range("E31").Value = "  -2.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.15"
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.34"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "55.44"
$ws.Range("E34").Value = "  -4.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.45"
$ws.Range("E35").Value = "  +1.95%  "
$ws.Range("E36").Value = "  -0.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "458.91"
$ws.Range("E37").Value = "  -7.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.232.39"
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0797"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0384"
$ws.Range("E40").Value = "  -3.51%  "
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.47"
$ws.Range("E43").Value = "  -7.03%  "
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.94"
$ws.Range("E45").Value = "  +5.76%  "
$ws.Range("E46").Value = "  -3.37%  "
$ws.Range("E47").Value = "  -1.83%  "
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.23"
$ws.Range("E49").Value = "  -3.84%  "
$ws.Range("E50").Value = "  -6.16%  "
$ws.Range("E51").Value = "  +7.26%  "
